$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You are a marketing researcher for a leading consumer goods company. You have been assigned to conduct an in-depth interview with a high-profile, notoriously reserved industry expert. You hope to gain valuable insights about emerging market trends that could impact your company's strategic plans. You will only have one interview with this expert, so you need to quickly build rapport and trust to encourage open and honest communication.What is the most effective approach to ensure a successful interview?",
        "ques_type": 2,
        "options": [
            "Begin the interview by discussing shared interests, then gradually transition into the main topic.",
            "Begin the interview by discussing your company's achievements, then gradually transition into the main topic. ",
            "Begin the interview by sharing your own opinions on emerging market trends, then gradually transition into the main topic.",
            "Begin the interview by asking the interviewee about their personal life, then gradually transition into the main topic."
        ],
        "score": "Begin the interview by discussing shared interests, then gradually transition into the main topic."
    },
    {
        "title": "As a content strategist for a tech startup, you are preparing to interview a subject matter expert (SME) in artificial intelligence (AI). You want to gather valuable insights and quotes for an upcoming blog post on the future of AI in healthcare. However, the SME is known for being highly technical in their explanations, which could potentially confuse your blog's target audience of non-technical healthcare professionals. What is the most effective strategy to elicit non-technical responses from the SME?",
        "ques_type": 2,
        "options": [
            "Request the SME provide a written response to your questions beforehand, allowing you to simplify their responses for your audience.",
            "Start the interview by explaining your audience's lack of technical knowledge, encouraging the SME to keep their explanations simple.",
            "Prepare a list of questions that encourages the SME to explain concepts in layperson\u2019s terms and refocus them if they become too technical.",
            "Ask the SME to provide a comprehensive list of all the technical terms they plan to use during the interview, allowing you to better prepare."
        ],
        "score": "Prepare a list of questions that encourages the SME to explain concepts in layperson\u2019s terms and refocus them if they become too technical."
    },
    {
        "title": "You are a lead researcher conducting an interview with a prominent CEO for an upcoming case study. The CEO prefers direct, to-the-point communication due to time constraints. Your goal is to extract as much relevant information as possible in the limited time you have without causing any discomfort or annoyance to the CEO. You decide to use both open-ended and closed-ended questions in the interview.What is the best approach to using closed-ended questions effectively in this scenario?",
        "ques_type": 2,
        "options": [
            "Use them to probe into the CEO's personal experiences and opinions.",
            "Use them to challenge the CEO\u2019s perspectives and elicit thought-provoking responses.\u00a0",
            "Use them to agree on interview logistics, such as how long the interview will be and whether you can record it.",
            "Use them to confirm details and extract concise responses on specific topics."
        ],
        "score": "Use them to confirm details and extract concise responses on specific topics."
    },
    {
        "title": "You are a hiring manager conducting an interview with a promising candidate for a key role in your company. The candidate\u2019s communication style is surprisingly informal, and it\u2019s hard to engage with them on a professional level.  How should you effectively communicate with the candidate?",
        "ques_type": 2,
        "options": [
            "Ignore the candidate's communication style and focus solely on their qualifications and experience.",
            "Request that the candidate adopt a more formal communication style for the remainder of the interview.",
            "Ask the candidate to explain why they choose to use such an informal communication style.",
            "Mirror the candidate's communication style while maintaining your professionalism."
        ],
        "score": "Mirror the candidate's communication style while maintaining your professionalism."
    }
]
'@

# The here-string adds a trailing newline; strip it so the cell content
# matches the source exactly (no trailing newline).
$newText = $newText.TrimEnd("`r", "`n")

# Remove the old second row entirely and put the full text (with its
# original bold/bordered style cleared) into A1.
$ws.Rows.Item(2).Delete()

$a1 = $ws.Range("A1")
$a1.Value = $newText
$a1.Style = "Normal"

# Writing multi-line text auto-grows the row height; restore the sheet's
# default (no explicit/custom row height), matching the source file.
$ws.Rows.Item(1).AutoFit()
